$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.145.69"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -4.51%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.080.32"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -4.89%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "543.84"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -6.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.60"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -10.86%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.078.15"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -4.71%  "

$ws.Range("E9").Value = "  -4.24%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.156"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -4.34%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.22"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -11.92%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.461"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -5.49%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000227"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.61%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.68"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -8.80%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.550.91"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -5.51%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.146.96"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -4.63%  "

$ws.Range("E17").Value = "  -2.97%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.082.54"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -4.84%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.62"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -6.89%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "487.33"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -9.66%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.32"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -8.25%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.704"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -5.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.20"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -7.47%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.50"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.96%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.12"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -9.96%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.11%  "

$ws.Range("E27").Value = "  -8.07%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.16"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -11.26%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.07%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.91"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -15.13%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.07"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -5.73%  "

$ws.Range("E32").Value = "  -6.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "58.83"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +7.53%  "

$ws.Range("E34").Value = "  -11.15%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.99"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -5.45%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.17"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -7.80%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "464.44"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -17.28%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.136.17"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.43%  "

$ws.Range("E39").Value = "  -13.46%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0794"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -7.46%  "

$ws.Range("E41").Value = "  -10.29%  "

$ws.Range("E42").Value = "  -5.36%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.55"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -11.28%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.251"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -12.44%  "

$ws.Range("E45").Value = "  +0.03%  "

$ws.Range("E46").Value = "  -12.06%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.57"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -6.82%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "119.60"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.67%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.107"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -4.25%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₃0513"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -7.17%  "

$ws.Range("E51").Value = "  -8.60%  "
